$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-5.779280513524995, -3.271673738956447, 4.346601516008363),
    @(3.140387788414934, 0.259726375341407, -4.785581156611421),
    @(-1.632258296012878, 0.6425724923610687, -3.22618693113327),
    @(-3.271841421723368, 0.07577018067240526, -1.077775649726385),
    @(-3.740465611219407, -0.2502757757902145, -5.230584308505059),
    @(-1.173786669969556, -1.206141140311958, -5.999948702752588),
    @(-0.527452439069747, -1.933494433760643, -1.995455801486972),
    @(-2.818732134997842, -1.367858927696945, 0.861171409487724),
    @(-4.338251754641532, -0.358771674335002, 0.811524987220763),
    @(-5.003720842301852, -0.2583636995404971, 1.431181490421301)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
